$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("second_canada_wave_summary")

# Insert a new column before column A, shifting "criteria"/"n_protests" to B/C.
$ws.Range("A1").EntireColumn.Insert()

# New "wave" header + values for the inserted column A.
$ws.Range("A1").Value = "wave"
$ws.Range("A2").Value = "Second Canada Wave (Labor)"
$ws.Range("A3").Value = "Second Canada Wave (Labor)"

# Match the header formatting used by B1/C1 (bold, centered).
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108  # xlCenter

# New row 4: Quebec Economy/inequality Protests 2014/2015, 47
$ws.Range("A4").Value = "Second Canada Wave (Labor)"
$ws.Range("B4").Value = "Quebec Economy/inequality Protests 2014/2015"
$ws.Range("C4").Value = 47
